# Edit: add two new weekly records for "Ajo" (Hortaliza, Vega Central
# Mapocho de Santiago) at the top of the existing history block (rows
# 188-269), pushing the prior rows down to 190-271.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 188:189; everything that was at 188 onward
# (through 269) shifts down to 190 onward (through 271).
$ws.Rows("188:189").Insert()

# New row 188: "Extra" quality record, 2022-09-21 (serial 44825).
$ws.Cells.Item(188, 1).Value  = 9
$ws.Cells.Item(188, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(188, 3).Value  = "Metropolitana"
$ws.Cells.Item(188, 4).Value  = 44825
$ws.Cells.Item(188, 5).Value  = 13
$ws.Cells.Item(188, 6).Value  = 100112003
$ws.Cells.Item(188, 7).Value  = "Ajo"
$ws.Cells.Item(188, 8).Value  = "Chino"
$ws.Cells.Item(188, 9).Value  = "Extra"
$ws.Cells.Item(188, 10).Value = 50
$ws.Cells.Item(188, 11).Value = 25000
$ws.Cells.Item(188, 12).Value = 25000
$ws.Cells.Item(188, 13).Value = 25000
$ws.Cells.Item(188, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(188, 15).Value = "China"
$ws.Cells.Item(188, 16).Value = 2500
$ws.Cells.Item(188, 17).Value = 10
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# New row 189: "Primera" quality record, same date.
$ws.Cells.Item(189, 1).Value  = 9
$ws.Cells.Item(189, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(189, 3).Value  = "Metropolitana"
$ws.Cells.Item(189, 4).Value  = 44825
$ws.Cells.Item(189, 5).Value  = 13
$ws.Cells.Item(189, 6).Value  = 100112003
$ws.Cells.Item(189, 7).Value  = "Ajo"
$ws.Cells.Item(189, 8).Value  = "Chino"
$ws.Cells.Item(189, 9).Value  = "Primera"
$ws.Cells.Item(189, 10).Value = 80
$ws.Cells.Item(189, 11).Value = 20000
$ws.Cells.Item(189, 12).Value = 20000
$ws.Cells.Item(189, 13).Value = 20000
$ws.Cells.Item(189, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(189, 15).Value = "China"
$ws.Cells.Item(189, 16).Value = 2000
$ws.Cells.Item(189, 17).Value = 10
$ws.Cells.Item(189, 18).Value = "Hortaliza"
